$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Shay")

# Add the new bird record into row 5 (next empty row)
$ws.Cells.Item(5, 7).Value2  = "45327"
$ws.Cells.Item(5, 8).Value2  = "European Gouldian"
$ws.Cells.Item(5, 9).Value2  = "East Europe"
$ws.Cells.Item(5, 10).Value2 = "14/04/2023"
$ws.Cells.Item(5, 11).Value2 = "Male"
$ws.Cells.Item(5, 12).Value2 = "A5342G"
$ws.Cells.Item(5, 13).Value2 = "97123"
$ws.Cells.Item(5, 14).Value2 = "98722"
$ws.Cells.Item(5, 15).Value2 = "yes"

$sortRange = $ws.Range("G2:O5")
$keyRange = $ws.Range("G2:G5")
$sortRange.Sort($keyRange, 1)

for ($r = 1; $r -le 5; $r++) {
    $line = "Row $r : "
    for ($c = 7; $c -le 15; $c++) {
        $v = $ws.Cells.Item($r, $c).Value2
        $line += "[$v]"
    }
    Write-Host $line
}
